$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update / overwrite existing cell values -------------------------------
$ws.Range("B2").Value = "CD"
$ws.Range("G4").Value = "CRU"
$ws.Range("C6").Value = "C"
$ws.Range("E8").Value = "C"
$ws.Range("G13").Value = "C/RU"

# --- Clear cells that become empty (remove the cell entirely) --------------
$ws.Range("B4").Clear()
$ws.Range("F6").Clear()
$ws.Range("F9").Clear()

# --- New cell values ---------------------------------------------------------
$ws.Range("C10").Value = "R"

$ws.Range("J4").Value = "Bearbeiten auch C -> aufteilen auf zwei Dienste möglich"
$ws.Range("J7").Value = "Bearbeiten auch C -> aufteilen auf zwei Schichten möglich"
$ws.Range("J10").Value = "Kontext1: Tag Assi zuordnen // Kontext2: Tag Schicht zuordnen"
$ws.Range("J13").Value = "Kontext1: Assi-Zuordnung zu Dienst entspricht ""update Dienst"" // Kontext2: Assi-Zuordnung zu Schicht entspricht ""create Dienst"""

$ws.Range("F26").Value = "cx"

# --- Column widths / row height ---------------------------------------------
# Column G: target stored width 7 (no bestFit)
$ws.Columns.Item(7).ColumnWidth = 6.083334
# Column J: new column, target stored width ~63.88671875
$ws.Columns.Item(10).ColumnWidth = 62.916667

# Row 13 grows taller to fit the new wrapped note in J13
$ws.Rows.Item(13).RowHeight = 28.8

# --- Selection ---------------------------------------------------------------
$null = $ws.Range("B21").Select()
